# Apply the corrected "jump" formulas to the weather sheet.
# The original formulas used IF(ABS(curr-prev) < 2, 0, 1); the corrected
# versions drop the ABS() and simply compare (prev-curr) < 2, i.e. they
# swap the operand order and remove the absolute value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: AA2:AW2 compare the current column to the one before it (no ABS, prev-curr).
# Z2 itself is a plain value cell (not a formula) and is left untouched.
$ws.Range("AA2:AW2").FormulaR1C1 = "=IF(RC[-25]-RC[-24]<2,0,1)"

# Column Z (rows 3:92): compares column Y of the previous row to column B of
# the current row.
$ws.Range("Z3:Z92").FormulaR1C1 = "=IF(R[-1]C[-1]-RC[-24]<2,0,1)"

# Columns AA:AW (rows 3:92): compares each column to the one immediately to
# its left, within the same row.
$ws.Range("AA3:AW92").FormulaR1C1 = "=IF(RC[-25]-RC[-24]<2,0,1)"

# Update the view: select Z3:Z92 with Z3 as the active cell (as in the edited file).
$ws.Range("Z3:Z92").Select()
